$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "Censu" -> "Census" in the annotation field description (row 7, column D)
$ws.Range("D7").Value = "Census notes on statistial exceptions that occurred for a given year, variable and geography"

# Update the active selection to match the authored state (D7)
$ws.Range("D7").Select()
